$wb = $excel.ActiveWorkbook

# --- Sheet: Forecast Comparison ---
$ws1 = $wb.Worksheets.Item("Forecast Comparison")
$ws1.Range("D2").Value = 197
$ws1.Range("D3").Value = 196
$ws1.Range("D10").Value = 161
$ws1.Range("D11").Value = 156
$ws1.Range("D15").Value = 138
$ws1.Range("D16").Value = 125
$ws1.Range("D17").Value = 114

# --- Sheet: Summary ---
# Column B on this sheet stores every value as plain text (inline strings in
# the source file), including numbers and dates. Force the NumberFormat to
# Text ("@") before assigning so Excel's COM layer doesn't auto-coerce the
# numeric-looking / date-looking strings into real numbers/dates.
$ws2 = $wb.Worksheets.Item("Summary")

$ws2.Range("B9").NumberFormat = "@"
$ws2.Range("B9").Value = "2505"

$ws2.Range("B10").NumberFormat = "@"
$ws2.Range("B10").Value = "1358"

$ws2.Range("B11").NumberFormat = "@"
$ws2.Range("B11").Value = "712"

$ws2.Range("B12").NumberFormat = "@"
$ws2.Range("B12").Value = "197"

$ws2.Range("B13").NumberFormat = "@"
$ws2.Range("B13").Value = "2025-01-26"

$ws2.Range("B14").NumberFormat = "@"
$ws2.Range("B14").Value = "114"
